$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cells (report volume/number + date range) ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Crime statistics table (rows 14-30) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("F14").NumberFormat = "general"
$ws.Range("H14").Value = -100
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = -55.555555555555
$ws.Range("L15").Value = -33.333333333333
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -11.764705882352
$ws.Range("I16").Value = 103
$ws.Range("J16").Value = 131
$ws.Range("K16").Value = -21.374045801526
$ws.Range("L16").Value = 10.752688172043
$ws.Range("M16").Value = 39.189189189189
$ws.Range("N16").Value = -77.657266811279
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 53.846153846153
$ws.Range("I17").Value = 118
$ws.Range("J17").Value = 95
$ws.Range("K17").Value = 24.210526315789
$ws.Range("L17").Value = 18
$ws.Range("M17").Value = 131.372549019608
$ws.Range("N17").Value = -26.708074534161
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -45.454545454545
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 36
$ws.Range("H18").Value = -47.222222222222
$ws.Range("I18").Value = 169
$ws.Range("J18").Value = 248
$ws.Range("K18").Value = -31.854838709677
$ws.Range("L18").Value = 49.557522123893
$ws.Range("M18").Value = 56.481481481481
$ws.Range("N18").Value = -62.444444444444
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 40
$ws.Range("E19").Value = -47.5
$ws.Range("F19").Value = 101
$ws.Range("G19").Value = 150
$ws.Range("H19").Value = -32.666666666666
$ws.Range("I19").Value = 683
$ws.Range("J19").Value = 717
$ws.Range("K19").Value = -4.741980474198
$ws.Range("L19").Value = 99.125364431486
$ws.Range("M19").Value = 16.952054794520
$ws.Range("N19").Value = -48.685199098422
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("I20").Value = 22
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -35.294117647058
$ws.Range("L20").Value = 4.761904761904
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -94.581280788177
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -39.655172413793
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 225
$ws.Range("H21").Value = -28.444444444444
$ws.Range("I21").Value = 1100
$ws.Range("J21").Value = 1236
$ws.Range("K21").Value = -11.003236245954
$ws.Range("L21").Value = 62.721893491124
$ws.Range("M21").Value = 30.177514792899
$ws.Range("N21").Value = -60.9375
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "general"
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 24
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = -27.272727272727
$ws.Range("C24").Value = 57
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = 42.5
$ws.Range("F24").Value = 223
$ws.Range("G24").Value = 165
$ws.Range("H24").Value = 35.151515151515
$ws.Range("I24").Value = 1166
$ws.Range("J24").Value = 1085
$ws.Range("K24").Value = 7.465437788018
$ws.Range("L24").Value = 71.218795888399
$ws.Range("M24").Value = 35.739231664726
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 8
$ws.Range("F25").Value = 39
$ws.Range("H25").Value = 18.181818181818
$ws.Range("I25").Value = 254
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 15.454545454545
$ws.Range("L25").Value = 70.469798657718
$ws.Range("M25").Value = 88.148148148148
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E26").Value = -100
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = -53.846153846153
$ws.Range("L26").Value = -25
$ws.Range("D27").Value = 5
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -36.363636363636
$ws.Range("J27").Value = 41
$ws.Range("K27").Value = -4.878048780487
$ws.Range("L27").Value = 34.482758620689
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D30").NumberFormat = "general"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E30").NumberFormat = "general"
$ws.Range("G30").Value = 1
